# Add Denmark, Sweden and Norway market test data sheets by duplicating the
# existing Belgium sheet (same template/layout), trimming the panel rows
# down to FC602S/FC604S only, and filling in the market-specific values.

$wb = $excel.ActiveWorkbook
$belgium = $wb.Worksheets.Item("Belgium")

# --- Denmark ---------------------------------------------------------
$belgium.Copy($null, $belgium)
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("A10:A13").EntireRow.Delete()
$denmark.Range("B2").Value = "Denmark market"
$denmark.Range("B4").Value = "NGC-3446/T2003"

# --- Sweden ------------------------------------------------------------
$belgium.Copy($null, $denmark)
$sweden = $wb.Worksheets.Item($wb.Worksheets.Count)
$sweden.Name = "Sweden"
$sweden.Range("A10:A13").EntireRow.Delete()
$sweden.Range("B2").Value = "Sweden market"
$sweden.Range("B4").Value = "NGC-3465/T2029"

# --- Norway --------------------------------------------------------------
$belgium.Copy($null, $sweden)
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("A10:A13").EntireRow.Delete()
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1918"

# Leave the Belgium / Denmark / Sweden sheets with the whole sheet selected
# (as left behind by the copy operations) and finish on Norway with cell B5
# selected, which is the active sheet/tab when the workbook is saved.
$belgium.Activate()
$belgium.Cells.Select()

$denmark.Activate()
$denmark.Cells.Select()

$sweden.Activate()
$sweden.Cells.Select()

$norway.Activate()
$norway.Range("B5").Select()
